# "updated main GSC export data"
#
# The workbook is a rolling-window GSC export:
#   - Sheet "Chart" (sheet1) has a header row followed by one row per day
#     (col A = date, col B = Non-HTTPS URLs, col C = HTTPS URLs).
#   - Sheet "Table" (sheet2) just has a 3-column header row ("Issue",
#     "Validation", "Pages") used elsewhere in the workbook.
#
# This edit rolls the window forward by one day:
#   - drop the oldest date (2025-11-14)
#   - every remaining day's data shifts up one row
#   - two new days are appended at the bottom (2026-02-12, 2026-02-13)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Chart"

# Drop the oldest day (row 2, 2025-11-14). Deleting the whole row shifts
# every later row (dates in col A together with their col B/C figures) up
# by one, which is exactly the "rolling window" shift seen across the rest
# of the sheet.
$ws.Rows.Item(2).Delete()

# After the delete, row 91 still holds the last existing day (2026-02-11).
# Append the two new days below it as plain text dates (not auto-converted
# date serials) with their figures.
$ws.Range("A91").NumberFormat = "@"
$ws.Range("A91").Value = "2026-02-12"
$ws.Range("A91").ClearFormats()
$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 30

$ws.Range("A92").NumberFormat = "@"
$ws.Range("A92").Value = "2026-02-13"
$ws.Range("A92").ClearFormats()
$ws.Range("B92").Value = 0
$ws.Range("C92").Value = 30
